$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new task rows (set column A text first for both rows so the
# shared-string table is built in the same order as the authored workbook)
$ws.Range("A25").Value = "New objects to define wining combinations"
$ws.Range("A26").Value = "New object to compare 2 results and determine wining combination"
$ws.Range("B25").Value = "ResultComparer"
$ws.Range("B26").Value = "ResultComparer"
$ws.Range("C25").Value = 1
$ws.Range("C26").Value = 1

# Update the selected cell to match the new active selection
$ws.Range("C25").Select()
